$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.629.40'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.634.08'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.31'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.33'
$ws.Range('E8').Value = '  +1.56%  '
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0863'
$ws.Range('E11').Value = '  -3.26%  '
$ws.Range('D12').Value = '1.864.58'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = '1.638.00'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.21'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').Value = '27.620.50'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.38'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('E20').Value = '  -2.33%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  +5.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.37'
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.12'
$ws.Range('E24').Value = '  +5.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.08'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E26').Value = '  -0.97%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.55'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('D33').Value = '1.475.45'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.09'
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.54'
$ws.Range('E35').Value = '  -2.02%  '
$ws.Range('E36').Value = '  -1.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.957'
$ws.Range('E37').Value = '  +7.00%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.880'
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.560'
$ws.Range('E39').Value = '  -0.57%  '
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('E41').Value = '  +2.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '67.85'
$ws.Range('E43').Value = '  -3.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.21'
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('E46').Value = '  -4.86%  '
$ws.Range('D47').Value = '1.774.39'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.70'
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('E51').Value = '  +0.55%  '
